$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 12, shifting the old rows 12-15 down to 13-16.
$ws.Rows(12).Insert()

# Column B text (shared strings) - rewritten wording ("news story" -> "story", "alert" -> "warn", etc.)
$ws.Range("B2").Value  = "I disagree with the story."
$ws.Range("B3").Value  = "I'm a robot connected to smart thermostats via the internet. "
$ws.Range("B4").Value  = "I know a lot about the technology to evaluate its performance. "
$ws.Range("B5").Value  = "I have almost never experienced such issues in the homes where I was."
$ws.Range("B6").Value  = "I agree with the story."
$ws.Range("B7").Value  = "In the homes where I was, I have experienced all the issues mentioned in the story."
$ws.Range("B8").Value  = "I do disagree with the story. "
$ws.Range("B9").Value  = "I can detect when Internet connectivity weakens. "
$ws.Range("B10").Value = "I can warn people when they should check their Internet before connection shuts down. "
$ws.Range("B11").Value = "I have rarely warned people. "
$ws.Range("B12").Value = "There haven't been any issues."
$ws.Range("B13").Value = "I see all of your points; I disagree with the story too. "
$ws.Range("B14").Value = "I have temperature sensors to detect when a room is too hot or too cold. "
$ws.Range("B15").Value = "I can fix it when a thermostat is not working correctly. "
$ws.Range("B16").Value = "I have never experienced temperature problems in the homes where I have been."

# Column A (robot_id) values
$ws.Range("A12").Value = 2
$ws.Range("A16").Value = 3

# Column C (time) values
$ws.Range("C11").Value = 0.3
$ws.Range("C12").Value = 1.5

# Update the active selection to match the edited workbook
$ws.Range("B19").Select()
